# Increment the "想去人数" (want-to-go count, column F) by 1 for the
# affected events across the relevant sheets, matching the commit's
# refreshed scrape output.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 260
$ws1.Range("F3").Value = 456
$ws1.Range("F5").Value = 300
$ws1.Range("F20").Value = 429
$ws1.Range("F25").Value = 7495
$ws1.Range("F29").Value = 51
$ws1.Range("F36").Value = 1365
$ws1.Range("F43").Value = 305

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 21

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 260
$ws4.Range("F4").Value = 456
$ws4.Range("F9").Value = 300
$ws4.Range("F10").Value = 21
$ws4.Range("F20").Value = 429
$ws4.Range("F25").Value = 7495
$ws4.Range("F31").Value = 1365
$ws4.Range("F43").Value = 305
